# Update the "Contabilidade" sheet so it reports three years (2020-2022)
# of financial-ratio data instead of a single column of ratio labels/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table layout (row 1 = header with the three fiscal-year labels,
# rows 2-9 = one financial ratio per row, one column per year).
#   A: ratio name            B: 31/2020   C: 31/2021   D: 31/2022
$data = @{
    "A1" = "Ano";                           "B1" = "31/2020"; "C1" = "31/2021"; "D1" = "31/2022"
    "A2" = "Líquidez Imediata";             "B2" = "0.44";    "C2" = "0.43";    "D2" = "0.25"
    "A3" = "Líquidez Corrente";             "B3" = "1.04";    "C3" = "1.24";    "D3" = "0.99"
    "A4" = "Líquidez Seca";                 "B4" = "0.83";    "C4" = "0.94";    "D4" = "0.71"
    "A5" = "Líquidez Geral";                "B5" = "0.30";    "C5" = "0.37";    "D5" = "0.44"
    "A6" = "Composição do Endividamento";   "B6" = "20.15%";  "C6" = "23.12%";  "D6" = "26.73%"
    "A7" = "PCT";                           "B7" = "80.0%";   "C7" = "59.95%";  "D7" = "62.69%"
    "A8" = "Margem Líquida";                "B8" = "2.29%";   "C8" = "23.69%";  "D8" = "29.47%"
    "A9" = "Rentabilidade";                 "B9" = "0.73%";   "C9" = "11.02%";  "D9" = "19.35%"
}

# Cells whose new text would otherwise be auto-recognized as a number/
# percentage by Excel's normal cell-entry parsing. The source workbook
# stores every one of these as literal text (shared string), so mark the
# cells as Text before typing the value, then drop the number format again
# so the cell keeps using the sheet's default (General) style.
$textCells = @(
    "B2","C2","D2",
    "B3","C3","D3",
    "B4","C4","D4",
    "B5","C5","D5",
    "B6","C6","D6",
    "B7","C7","D7",
    "B8","C8","D8",
    "B9","C9","D9"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}

# Reset the saved selection back to the default top-left cell.
$ws.Range("A1").Select()
